# Auto-generated edit script: updates the Price (D) and Volume(1h) (E)
# columns of the cryptos worksheet to reflect refreshed market data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.751.47"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "'2.766.00"
$ws.Range("E3").Value = "  -2.47%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'354.69"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").Value = "'109.25"
$ws.Range("E6").Value = "  -4.32%  "
$ws.Range("D7").Value = "'0.563"
$ws.Range("E7").Value = "  +1.64%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.589"
$ws.Range("E9").Value = "  -2.19%  "
$ws.Range("D10").Value = "'39.95"
$ws.Range("E10").Value = "  -4.08%  "
$ws.Range("D11").Value = "'0.0849"
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").Value = "'19.33"
$ws.Range("E13").Value = "  -3.64%  "
$ws.Range("D14").Value = "'7.58"
$ws.Range("E14").Value = "  -2.15%  "
$ws.Range("D15").Value = "'3.208.72"
$ws.Range("E15").Value = "  -2.05%  "
$ws.Range("D16").Value = "'2.783.84"
$ws.Range("E16").Value = "  -2.00%  "
$ws.Range("D17").Value = "'0.928"
$ws.Range("E17").Value = "  +3.10%  "
$ws.Range("D18").Value = "'51.659.80"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("D19").Value = "'7.41"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "'3.14"
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("D21").Value = "'13.00"
$ws.Range("E21").Value = "  -4.23%  "
$ws.Range("D22").Value = "'0.0₃0973"
$ws.Range("E22").Value = "  -2.82%  "
$ws.Range("D23").Value = "'273.88"
$ws.Range("E23").Value = "  +0.94%  "
$ws.Range("D24").Value = "'69.77"
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("D25").Value = "'2.71"
$ws.Range("E25").Value = "  -2.95%  "
$ws.Range("D26").Value = "'26.42"
$ws.Range("E26").Value = "  -1.45%  "
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").Value = "'10.11"
$ws.Range("E28").Value = "  -2.33%  "
$ws.Range("D29").Value = "'2.21"
$ws.Range("E29").Value = "  -1.70%  "
$ws.Range("E30").Value = "  +2.31%  "
$ws.Range("D31").Value = "'51.57"
$ws.Range("E31").Value = "  +1.62%  "
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("D33").Value = "'33.75"
$ws.Range("E33").Value = "  -0.68%  "
$ws.Range("D34").Value = "'5.68"
$ws.Range("E34").Value = "  -2.83%  "
$ws.Range("D35").Value = "'0.0841"
$ws.Range("E35").Value = "  +0.99%  "
$ws.Range("E36").Value = "  +5.93%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("E38").Value = "  -0.43%  "
$ws.Range("D39").Value = "'17.98"
$ws.Range("E39").Value = "  -0.85%  "
$ws.Range("E40").Value = "  -5.09%  "
$ws.Range("E41").Value = "  -2.33%  "
$ws.Range("E42").Value = "  -0.97%  "
$ws.Range("E43").Value = "  -2.91%  "
$ws.Range("D44").Value = "'120.85"
$ws.Range("E44").Value = "  -4.40%  "
$ws.Range("D45").Value = "'21.85"
$ws.Range("E45").Value = "  -8.35%  "
$ws.Range("D46").Value = "'2.061.26"
$ws.Range("E46").Value = "  -1.50%  "
$ws.Range("D47").Value = "'3.22"
$ws.Range("E47").Value = "  -4.51%  "
$ws.Range("D48").Value = "'2.21"
$ws.Range("E48").Value = "  -4.35%  "
$ws.Range("E49").Value = "  -1.27%  "
$ws.Range("E50").Value = "  -2.75%  "
$ws.Range("D51").Value = "'8.92"
$ws.Range("E51").Value = "  -0.47%  "
